$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Category" to A1, matching the style of the other header cells (B1:W1):
# bold font, centered horizontally, top-aligned vertically, thin border all around.
$ws.Range("A1").Value = "Category"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4160
$ws.Range("A1").Borders.LineStyle = 1

# Remove the bold/border style from A2:A46 (they previously used the same style as the header)
$ws.Range("A2:A46").Style = "Normal"
